# Adds a new "Form Tag" column (V) to the CapitalCommitment sheet.
# Header V1 = "Form Tag" (mirrors the other header cells in row 1).
# Data cells V2:V9 = "Default" (mirrors the other data cells in their rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy formatting from the last existing header cell (U1)
# so the new column blends in with the rest of the header row.
$ws.Range("V1").Value = "Form Tag"
$ws.Range("V1").Style = $ws.Range("U1").Style

# Data cells - copy formatting from the corresponding row's existing
# "CF 3" cell (column U) so each new cell matches its row's look.
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 22)
    $cell.Value = "Default"
    $cell.Style = $ws.Cells.Item($r, 21).Style
}

# Match the cursor/selection position left behind by the edit session.
[void]$ws.Range("W11").Select()
